$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Session 4 (Greedy) mark for the student, row 4, column E
$ws.Range("E4").Value = 8

# Feedback comment for Session 4 (Greedy), row 5, column E
$ws.Range("E5").Value = "Good but the expected complexity for greedy 2 and greedy 3 is O(nlogn) if you sort the elements beforehand or if you use a priority queue"

# Update the selection to match the newly-filled range, like the author did after entering data
$ws.Range("E5:E12").Select()
